# Fruta / hortaliza, semanal
# Insert a new weekly price-observation row for "Alcachofa" (Española,
# Primera) at row 318, shifting the existing rows 318:406 down to 319:407.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 318 (pushes 318:406 -> 319:407)
$ws.Rows("318:318").Insert()

# Populate the newly inserted row with the new weekly record
$ws.Range("A318").Value = 9
$ws.Range("B318").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C318").Value = "Metropolitana"
$ws.Range("D318").Value = 44736
$ws.Range("E318").Value = 13
$ws.Range("F318").Value = 100112013
$ws.Range("G318").Value = "Alcachofa"
$ws.Range("H318").Value = "Española"
$ws.Range("I318").Value = "Primera"
$ws.Range("J318").Value = 52
$ws.Range("K318").Value = 22000
$ws.Range("L318").Value = 22000
$ws.Range("M318").Value = 22000
$ws.Range("N318").Value = "$/caja 30 unidades"
$ws.Range("O318").Value = "Provincia del Elquí"
$ws.Range("P318").Value = 733
$ws.Range("Q318").Value = 30
$ws.Range("R318").Value = "Hortaliza"

Write-Host "Inserted new row 318; sheet now spans through row 407."
